$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.841.75"
$ws.Range("E2").Value = "  +1.78%  "

$ws.Range("D3").Value = "3.456.48"
$ws.Range("E3").Value = "  +1.24%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'575.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.04%  "

$ws.Range("D6").Value = "'161.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.48%  "

$ws.Range("D7").Value = "'0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +13.72%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").Value = "3.458.16"
$ws.Range("E9").Value = "  +1.33%  "

$ws.Range("D10").Value = "'7.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.40%  "

$ws.Range("D11").Value = "'0.124"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.86%  "

$ws.Range("E12").Value = "  +3.30%  "

$ws.Range("D13").Value = "4.065.90"
$ws.Range("E13").Value = "  +1.66%  "

$ws.Range("E14").Value = "  +0.31%  "

$ws.Range("D15").Value = "'0.0000191"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.15%  "

$ws.Range("D16").Value = "'28.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.53%  "

$ws.Range("D17").Value = "64.977.00"
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("D18").Value = "3.505.58"
$ws.Range("E18").Value = "  +2.39%  "

$ws.Range("D19").Value = "'6.43"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.48%  "

$ws.Range("D20").Value = "'14.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.80%  "

$ws.Range("D21").Value = "'379.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.31%  "

$ws.Range("D22").Value = "'8.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "

$ws.Range("D23").Value = "'0.549"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").Value = "'72.40"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.65%  "

$ws.Range("D26").Value = "'0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.44%  "

$ws.Range("D27").Value = "'9.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.27%  "

$ws.Range("D28").Value = "'0.178"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").Value = "'1.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.34%  "

$ws.Range("D31").Value = "'6.12"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("D32").Value = "'2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.63%  "

$ws.Range("D33").Value = "'23.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.28%  "

$ws.Range("D34").Value = "'7.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.90%  "

$ws.Range("D35").Value = "'1.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.82%  "

$ws.Range("D36").Value = "'161.72"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").Value = "'1.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.70%  "

$ws.Range("D38").Value = "'0.0774"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.09%  "

$ws.Range("D39").Value = "2.964.48"
$ws.Range("E39").Value = "  -0.58%  "

$ws.Range("D40").Value = "'26.51"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.77%  "

$ws.Range("D41").Value = "'6.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.43%  "

$ws.Range("D42").Value = "'4.54"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.32%  "

$ws.Range("D43").Value = "'0.0319"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "

$ws.Range("D44").Value = "'42.84"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.04%  "

$ws.Range("D45").Value = "'0.776"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.95%  "

$ws.Range("D46").Value = "'25.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.14%  "

$ws.Range("E47").Value = "  +2.50%  "

$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "'0.111"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.14%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'312.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.92%  "

$ws.Range("D50").Value = "'6.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.40%  "

$ws.Range("D51").Value = "'0.864"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.27%  "
